$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (swap with former row 6 content)
$ws.Range("D4").Value = 44334
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("Q4").Value = "$/caja 12 kilos empedrada"
$ws.Range("S4").Value = 1042
$ws.Range("T4").Value = 12

# Row 5 (swap with former row 7 content)
$ws.Range("D5").Value = 44330
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 15000
$ws.Range("P5").Value = 15500
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 861
$ws.Range("T5").Value = 18

# Row 6 (swap with former row 4 content)
$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 17500
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17750
$ws.Range("Q6").Value = "$/caja 16 kilos granel"
$ws.Range("S6").Value = 1109
$ws.Range("T6").Value = 16

# Row 7 (swap with former row 5 content)
$ws.Range("D7").Value = 44316
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 16000
$ws.Range("P7").Value = 16000
$ws.Range("Q7").Value = "$/caja 16 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 16
